# Adds a new diary entry for 2022-06-07 (Tuesday) right after the
# "今天又是陶源大帅逼的一天" paragraph, and drops the stray paragraph-mark
# rFonts/eastAsia hint on that paragraph (matches the target OOXML diff).

$d = $word.ActiveDocument

# Locate the target paragraph by its text instead of a hard-coded index,
# so the script is resilient to any pre-existing paragraph-count drift.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "今天又是陶源大帅逼的一天") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find target paragraph"
}

$target = $d.Paragraphs.Item($targetIndex)
$targetRange = $target.Range

$body = '<w:p><w:pPr><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>今天又是陶源大帅逼的一天</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>022</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>年6月7日星期二</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t>呦吼 陶源大帅逼又来咯</w:t></w:r></w:p>'
$wordOpenXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes" ?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + "<w:body>" + $body + "</w:body>" + '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($wordOpenXml)

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"

